$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new monthly row (01-07-2021) mirroring the previous row's totals.
#
# A plain `Range.Value = "01-07-2021"` gets auto-interpreted as a date by
# Excel's smart typing (since the text looks like a date) and would be
# stored as a date serial with a new number-format style, which doesn't
# match the source workbook (where this column is plain shared-string
# text with no cell style). To reproduce the literal text entry, write it
# as a formula-literal first (never date-coerced), then copy/paste-special
# as values so the final cell holds a plain text value with no formula and
# no extra style, exactly like the existing cells in this column.
$ws.Range("A96").Formula = '="01-07-2021"'
$ws.Range("A96").Copy()
$ws.Range("A96").PasteSpecial(-4163)
$ws.Range("B96").Value = 202
$ws.Range("C96").Value = 0
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 202
$ws.Range("K96").Value = 0
